$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.691.31"
$ws.Range("E2").Value = "  -0.71%  "
$ws.Range("D3").Value = "1.890.34"
$ws.Range("E3").Value = "  -1.09%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.42%  "
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4888"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.28%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2929"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06687"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.54%  "
$ws.Range("D10").Value = "1.890.69"
$ws.Range("E10").Value = "  -1.04%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.75"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.58%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07241"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.19%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "89.46"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.07%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.015"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.91%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6661"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.42%  "
$ws.Range("D16").Value = "30.641.29"
$ws.Range("E16").Value = "  -0.81%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000007909"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.64%  "
$ws.Range("E18").Value = "  -0.05%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.05"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.93%  "
$ws.Range("D20").Value = "2.134.76"
$ws.Range("E20").Value = "  -0.86%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.751"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.82%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "192.39"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.28%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.079"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.78%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.319"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.63%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "159.89"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.53%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.33"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.23%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.834"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.403"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.66%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.273"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.99%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09034"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.38%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.948"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.00%  "
$ws.Range("E33").Value = "  -1.22%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7347"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.99%  "
$ws.Range("E35").Value = "  -4.49%  "
$ws.Range("E36").Value = "  +0.55%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01825"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.669"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9252"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.30%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.052"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.11%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4427"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.21%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "104.76"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9998"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("E44").Value = "  -1.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1342"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.39%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.369"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.38%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4184"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.23%  "
$ws.Range("E48").Value = "  -0.35%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.681"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.64%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.411"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "33.30"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.32%  "
